# Set the "Industries" column (H) values to 0 for rows 31 through 176.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H31:H176").Value = 0
